# modify the detailUserActivity for the user info
# Adds two new notes about dialog dismiss() usage to Sheet2, just below the
# existing "단순하게 이름과 전화번호만 필요하다." note (row 7), leaving row 8 blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("B9").Value  = "다이얼로그를 닫을 땐 dismiss()를 쓴다."
$ws.Range("B10").Value = "그리고 dissmiss될 때 "

$ws.Range("F15").Select()
